$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Fitness column (C) values for rows 2-123 to reflect the
# corrected simulated-annealing algorithm output ("correction in sa
# algorithm and 746 logs").
$ws.Range("C2:C8").Value = 7909
$ws.Range("C9:C26").Value = 7889
$ws.Range("C27:C35").Value = 7887
$ws.Range("C36:C37").Value = 7765
$ws.Range("C38:C38").Value = 7728
$ws.Range("C39:C80").Value = 7310
$ws.Range("C81:C123").Value = 7293
